$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-data-mapping rows appended to the bottom of the sheet.
# Shared-string table insertion order (verified against target diff) is:
#   YouthPlacementPreferences, PlacementNeeds, ParentGuardPlacementPreferences
# even though on the sheet they end up ordered (row171=ParentGuard, row172=Youth,
# row173=PlacementNeeds), so cell values are written in that same order below.

$ws.Range("A172").Value = "YouthPlacementPreferences"
$ws.Range("B172").Value = "cares\Placement.xlsx"
$ws.Range("C172").Value = "YouthPlacementPreferences"
$ws.Range("D172").Value = 1

$ws.Range("A173").Value = "PlacementNeeds"
$ws.Range("B173").Value = "cares\Placement.xlsx"
$ws.Range("C173").Value = "PlacementNeeds"
$ws.Range("D173").Value = 1

$ws.Range("A171").Value = "ParentGuardPlacementPreferences"
$ws.Range("B171").Value = "cares\Placement.xlsx"
$ws.Range("C171").Value = "ParentGuardPlacementPreferences"
$ws.Range("D171").Value = 1

# Column C needs to widen to fit the new, longer text.
$ws.Columns("C").AutoFit()

# Move the active selection to the new last row, matching where the author
# ended up after typing the new rows.
$ws.Range("A173").Select()

# Page was set up for portrait printing as part of this save.
$ws.PageSetup.Orientation = 1
